$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the existing row 704 (old row 704
# becomes row 706, and so on through the end of the sheet). Excel copies
# the formatting (incl. the date number format on column D) from the row
# above the insertion point, matching the target's style on the new rows.
$ws.Rows.Item(704).Insert()
$ws.Rows.Item(704).Insert()

# Populate the first new row (704).
$ws.Cells.Item(704, 1).Value = 9
$ws.Cells.Item(704, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(704, 3).Value = "Metropolitana"
$ws.Cells.Item(704, 4).Value = 44706
$ws.Cells.Item(704, 5).Value = 13
$ws.Cells.Item(704, 6).Value = "Fruta"
$ws.Cells.Item(704, 7).Value = 100102
$ws.Cells.Item(704, 8).Value = "Cítricos"
$ws.Cells.Item(704, 9).Value = 100102005
$ws.Cells.Item(704, 10).Value = "Naranja"
$ws.Cells.Item(704, 11).Value = "Fukumoto"
$ws.Cells.Item(704, 12).Value = "Primera"
$ws.Cells.Item(704, 13).Value = 380
$ws.Cells.Item(704, 14).Value = 9000
$ws.Cells.Item(704, 15).Value = 9000
$ws.Cells.Item(704, 16).Value = 9000
$ws.Cells.Item(704, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(704, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(704, 19).Value = 500
$ws.Cells.Item(704, 20).Value = 18

# Populate the second new row (705).
$ws.Cells.Item(705, 1).Value = 9
$ws.Cells.Item(705, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(705, 3).Value = "Metropolitana"
$ws.Cells.Item(705, 4).Value = 44706
$ws.Cells.Item(705, 5).Value = 13
$ws.Cells.Item(705, 6).Value = "Fruta"
$ws.Cells.Item(705, 7).Value = 100102
$ws.Cells.Item(705, 8).Value = "Cítricos"
$ws.Cells.Item(705, 9).Value = 100102005
$ws.Cells.Item(705, 10).Value = "Naranja"
$ws.Cells.Item(705, 11).Value = "Valencia"
$ws.Cells.Item(705, 12).Value = "Primera"
$ws.Cells.Item(705, 13).Value = 350
$ws.Cells.Item(705, 14).Value = 12500
$ws.Cells.Item(705, 15).Value = 12500
$ws.Cells.Item(705, 16).Value = 12500
$ws.Cells.Item(705, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(705, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(705, 19).Value = 694
$ws.Cells.Item(705, 20).Value = 18
